$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    -1.5707963267948966,
    -1.2217304763960306,
    -0.87266462599716477,
    -0.52359877559829882,
    -0.3490658503988659,
    -0.17453292519943295,
    -0.15707963267948966,
    -0.13962634015954636,
    -0.12217304763960307,
    -0.10471975511965978,
    -0.087266462599716474,
    -0.069813170079773182,
    -0.05235987755982989,
    -0.034906585039886591,
    -0.017453292519943295,
    0,
    0.017453292519943295,
    0.034906585039886591,
    0.05235987755982989,
    0.069813170079773182,
    0.087266462599716474,
    0.10471975511965978,
    0.12217304763960307,
    0.13962634015954636,
    0.15707963267948966,
    0.17453292519943295,
    0.3490658503988659,
    0.52359877559829882,
    0.87266462599716477,
    1.2217304763960306,
    1.5707963267948966
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

$ws.Range("G8").Select()
